$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1233.8966

$ws.Range("H62").Value = 2732.8
$ws.Range("J62").Value = 2888.3333
$ws.Range("L62").Value = 2888.3333
$ws.Range("N62").Value = -4136.3333

$ws.Range("H65").Value = 2732.8
$ws.Range("J65").Value = 2888.3333
$ws.Range("L65").Value = 14441.6665
$ws.Range("N65").Value = -20681.6665

$ws.Range("H92").Value = 594.8
$ws.Range("I92").Value = 683.53845
$ws.Range("J92").Value = 18
$ws.Range("K92").Value = 683.53845
$ws.Range("L92").Value = 18
$ws.Range("M92").Value = 564.46155
$ws.Range("N92").Value = -2514

$ws.Range("H96").Value = 192.4
$ws.Range("I96").Value = 21
$ws.Range("J96").Value = 449.5
$ws.Range("K96").Value = 63
$ws.Range("L96").Value = 1348.5
$ws.Range("M96").Value = 1310
$ws.Range("N96").Value = -4094.5

$ws.Range("H137").Value = 3072.7646
$ws.Range("I137").Value = 2414.1428
$ws.Range("K137").Value = 7242.428400000001
$ws.Range("M137").Value = -4692.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5971.0557
$ws.Range("I32").Value = 6234.0586
$ws.Range("J32").Value = 1500
$ws.Range("K32").Value = 6234.0586
$ws.Range("L32").Value = 1500
$ws.Range("M32").Value = -5947.0586
$ws.Range("N32").Value = -2074

$ws.Range("H97").Value = 1602
$ws.Range("I97").Value = 1602
$ws.Range("K97").Value = 1602
$ws.Range("M97").Value = -1106

$ws.Range("H102").Value = 2086.625
$ws.Range("I102").Value = 2241.8572
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 2241.8572
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = -619.8571999999999
$ws.Range("N102").Value = -4244

$ws.Range("H121").Value = 125000
$ws.Range("J121").Value = 125000
$ws.Range("L121").Value = 125000
$ws.Range("N121").Value = -128494

$ws.Range("H122").Value = 2487.8572
$ws.Range("I122").Value = 2880.2
$ws.Range("K122").Value = 8640.599999999999
$ws.Range("M122").Value = -6190.599999999999

$ws.Range("H133").Value = 140000
$ws.Range("J133").Value = 140000
$ws.Range("L133").Value = 140000
$ws.Range("N133").Value = -145060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 677
$ws.Range("I80").Value = 49.5
$ws.Range("K80").Value = 49.5
$ws.Range("M80").Value = 948.5

$ws.Range("H83").Value = 677
$ws.Range("I83").Value = 49.5
$ws.Range("K83").Value = 247.5
$ws.Range("M83").Value = 4744.5

$ws.Range("H99").Value = 3277.4
$ws.Range("I99").Value = 3096.75
$ws.Range("K99").Value = 3096.75
$ws.Range("M99").Value = -1598.75

$ws.Range("H134").Value = 1972.8572
$ws.Range("I134").Value = 1972.8572
$ws.Range("K134").Value = 5918.571599999999
$ws.Range("M134").Value = -3383.571599999999

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1107.25
$ws.Range("I5").Value = 715
$ws.Range("J5").Value = 1499.5
$ws.Range("K5").Value = 715
$ws.Range("L5").Value = 1499.5
$ws.Range("M5").Value = -603
$ws.Range("N5").Value = -1723.5

$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = ""

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").Value = ""

$ws.Range("H134").Value = 4006
$ws.Range("I134").Value = 4006
$ws.Range("K134").Value = 12018
$ws.Range("M134").Value = -9483

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 799.58826
$ws.Range("I5").Value = 616.5
$ws.Range("J5").Value = 1005.5625
$ws.Range("K5").Value = 1849.5
$ws.Range("L5").Value = 3016.6875
$ws.Range("M5").Value = -1737.5
$ws.Range("N5").Value = -3240.6875

$ws.Range("H92").Value = 499.66666
$ws.Range("I92").Value = 499
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 1497
$ws.Range("L92").Value = 1500
$ws.Range("M92").Value = -249
$ws.Range("N92").Value = -3996

$ws.Range("H132").Value = 3281
$ws.Range("J132").Value = 2993.4
$ws.Range("L132").Value = 26940.6
$ws.Range("N132").Value = -32000.6

$ws.Range("H135").Value = 799.58826
$ws.Range("I135").Value = 616.5
$ws.Range("J135").Value = 1005.5625
$ws.Range("K135").Value = 5548.5
$ws.Range("L135").Value = 9050.0625
$ws.Range("M135").Value = -3013.5
$ws.Range("N135").Value = -14120.0625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 514.875
$ws.Range("I102").Value = 514.875
$ws.Range("K102").Value = 514.875
$ws.Range("M102").Value = 1107.125

$ws.Range("H126").Value = 3999.5
$ws.Range("I126").Value = 3999.5
$ws.Range("K126").Value = 11998.5
$ws.Range("M126").Value = -9528.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3249.1667
$ws.Range("J46").Value = 3123.75
$ws.Range("L46").Value = 3123.75
$ws.Range("N46").Value = -3499.75

$ws.Range("H57").Value = 4000
$ws.Range("I57").Value = 4000
$ws.Range("K57").Value = 4000
$ws.Range("M57").Value = -3434

$ws.Range("H58").Value = 26678.6
$ws.Range("I58").Value = 5464.3335
$ws.Range("J58").Value = 58500
$ws.Range("K58").Value = 5464.3335
$ws.Range("L58").Value = 58500
$ws.Range("M58").Value = -5204.3335
$ws.Range("N58").Value = -59020

$ws.Range("H61").Value = 3299.8
$ws.Range("I61").Value = 3999.6667
$ws.Range("J61").Value = 2250
$ws.Range("K61").Value = 3999.6667
$ws.Range("L61").Value = 2250
$ws.Range("M61").Value = -3797.6667
$ws.Range("N61").Value = -2654

$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").Value = ""

$ws.Range("H100").Value = 3933
$ws.Range("I100").Value = 3999.5
$ws.Range("J100").Value = 3800
$ws.Range("K100").Value = 3999.5
$ws.Range("L100").Value = 3800
$ws.Range("M100").Value = -3458.5
$ws.Range("N100").Value = -4882

$ws.Range("H113").Value = 3299.8
$ws.Range("I113").Value = 3999.6667
$ws.Range("J113").Value = 2250
$ws.Range("K113").Value = 3999.6667
$ws.Range("L113").Value = 2250
$ws.Range("M113").Value = -1829.6667
$ws.Range("N113").Value = -6590

$ws.Range("H133").Value = 130000
$ws.Range("J133").Value = 130000
$ws.Range("L133").Value = 130000
$ws.Range("N133").Value = -135060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").Value = ""

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = ""

Write-Output "done"